# Update legacy GSC export data ("Chart" sheet).
#
# The rolling date window advances by one day:
#   - every date label in column A (rows 2..last) moves forward one day
#     (2025-10-17..2026-01-14 becomes 2025-10-18..2026-01-15)
#   - the HTTPS URL counts in column C shift up by one row to stay lined
#     up with the new date window; the freshly observed last day gets 0
#   - column B (Non-HTTPS URLs) is untouched (always 0)
#
# Column A holds the dates as literal text (not real Excel dates), so we
# shift it with Copy (cell -> cell) rather than typing the text via
# .Value: typing a date-shaped string into a General-formatted cell makes
# Excel auto-recognise it as a date and reformat the cell, which would
# introduce style churn that isn't part of this edit. Copying an existing
# text cell carries its literal string type/format across untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$lastRow = $ws.UsedRange.Rows.Count   # last data row (91: header + 90 days)

# --- Column C (HTTPS URLs): capture current values, then shift up by one row ---
$cValues = @{}
for ($r = 2; $r -le $lastRow; $r++) {
    $cValues[$r] = $ws.Cells.Item($r, 3).Value2
}
for ($r = 2; $r -lt $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = $cValues[$r + 1]
}
$ws.Cells.Item($lastRow, 3).Value = 0.0

# --- Column A (Date labels): shift every date forward by one day ---
# Row r takes on the text that currently lives in row r+1.
for ($r = 2; $r -lt $lastRow; $r++) {
    $ws.Cells.Item($r + 1, 1).Copy($ws.Cells.Item($r, 1))
}

# The new final day (one after the old last date) has no donor cell yet,
# so build it via a scratch formula (keeps it a plain string, avoiding
# date auto-detection) and copy the computed text into place.
$scratch = $ws.Cells.Item(1, 26)   # Z1 - unused scratch cell
$scratch.Formula = "=""2026-01-15"""
$scratch.Copy($ws.Cells.Item($lastRow, 1))
$scratch.ClearContents()
